$wb = $excel.ActiveWorkbook

# "Repayment schedule" sheet (sheet3.xml) gets a new blank column inserted
# before column N (14th column), shifting the old N/O/P (Late/Outstanding/
# Disbursement) columns one place to the right.
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

$wsRepay.Columns.Item(14).Insert() | Out-Null

# Match the width Excel gives the freshly inserted column (inherits the
# width of the column immediately to its left).
$wsRepay.Columns.Item(14).ColumnWidth = $wsRepay.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with the new selection.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("K14").Select() | Out-Null

# "Transactions" sheet (sheet4.xml) is no longer the selected tab; select a
# neutral range so it keeps its previous selection but loses tabSelected.
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2:L3").Select() | Out-Null

# Re-activate "Repayment schedule" so it stays the workbook's active sheet.
$wsRepay.Activate() | Out-Null
